# Add a new "部门" (Department) column between "职位" (Position) and
# "座位号" (Seat number) — i.e. insert a new column F, pushing the old
# F (座位号) and G (分机号) columns to G and H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh column at F; existing F/G (and their formatting) shift to G/H.
$ws.Columns.Item(6).Insert() | Out-Null

# Header for the new column.
$ws.Range("F1").Value = "部门"

# Department values for each team member row.
$ws.Range("F2").Value = "管理部"
$ws.Range("F3").Value = "设计部"
$ws.Range("F4").Value = "开发部"
$ws.Range("F5").Value = "客服部"
$ws.Range("F6").Value = "客服部"
$ws.Range("F7").Value = "客服部"

# Match the author's final selection.
$ws.Range("F14").Select() | Out-Null
